$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 08:38"

# Row 39 - Ucrania
$ws.Cells.Item(39, 2).Value = 54771
$ws.Cells.Item(39, 3).Value = 638
$ws.Cells.Item(39, 4).Value = 27154
$ws.Cells.Item(39, 5).Value = 26205
$ws.Cells.Item(39, 7).Value = 14
$ws.Cells.Item(39, 8).Value = 1412

# Row 48 - Afganistan
$ws.Cells.Item(48, 2).Value = 34730
$ws.Cells.Item(48, 3).Value = 275
$ws.Cells.Item(48, 4).Value = 21417
$ws.Cells.Item(48, 5).Value = 12275
$ws.Cells.Item(48, 7).Value = 26
$ws.Cells.Item(48, 8).Value = 1038

# Row 67 - Uzbekistan
$ws.Cells.Item(67, 2).Value = 13872
$ws.Cells.Item(67, 3).Value = 281
$ws.Cells.Item(67, 5).Value = 5778

# Row 76 - El Salvador
$ws.Cells.Item(76, 4).Value = 5761
$ws.Cells.Item(76, 5).Value = 3950

# Row 144 - Georgia
$ws.Cells.Item(144, 2).Value = 999
$ws.Cells.Item(144, 3).Value = 4
$ws.Cells.Item(144, 4).Value = 870
$ws.Cells.Item(144, 5).Value = 114
